$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bosses")

# Update the "Name" column (B) values for each boss row first.
$ws.Range("B2").Value = "BossA"
$ws.Range("B3").Value = "BossB"
$ws.Range("B4").Value = "BossC"
$ws.Range("B5").Value = "BossD"
$ws.Range("B6").Value = "BossE"
$ws.Range("B7").Value = "BossF"

# Then update the "Id" column (A) values for each boss row.
$ws.Range("A2").Value = "B_A"
$ws.Range("A3").Value = "B_B"
$ws.Range("A4").Value = "B_C"
$ws.Range("A5").Value = "B_D"
$ws.Range("A6").Value = "B_A"
$ws.Range("A7").Value = "B_F"

# Update the current selection on the sheet to match the author's position.
$ws.Range("B9").Select()
